$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.114.05"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.31%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.812.14"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.89%  "

$ws.Range("E4").Value = "  +0.10%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "311.86"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4607"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +4.80%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3742"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.97%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07397"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8600"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.59"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.820.41"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.83%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.653"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.84%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.374"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.41%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.07091"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.47%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "91.60"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("E17").Value = "  +0.24%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008741"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +1.04%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "27.123.30"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.15%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.318"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.63%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.87"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.044.15"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.78%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "1.928"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.20%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "151.80"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.203"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.72%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "18.43"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.263"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.78%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "116.66"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08894"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7712"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.11%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.167"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.18%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.519"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.06%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.888"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("E36").Value = "  +0.20%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.127"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.04%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01959"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.27%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05227"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.38%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "7.273"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.25%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.386"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +20.85%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.920"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.20%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.5278"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1679"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.49%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.587"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.14%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5029"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.55%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.32"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "104.88"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("E49").Value = "  +0.22%  "

$ws.Range("E50").Value = "  +0.85%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06317"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.21%  "
